# Finish the 10th scenario skeleton:
#  - rename ID_0010_NO_DATA_ENTITY_REQUIRED -> ID_0010
#  - fill its sheet with the same USUARIO/SENHA/CATEGORIA sample row used by
#    the neighbouring ID_0009 sheet
#  - move the "active" selection/tab over from ID_0008 to the now-finished
#    ID_0010 sheet, and nudge ID_0009's selection down to include the data row

$wb = $excel.ActiveWorkbook

$ws0008   = $wb.Worksheets.Item("ID_0008")
$ws0009   = $wb.Worksheets.Item("ID_0009")
$ws0010   = $wb.Worksheets.Item("ID_0010_NO_DATA_ENTITY_REQUIRED")

# 1. rename the sheet now that it is no longer "no data entity required"
$ws0010.Name = "ID_0010"

# 2. populate ID_0010 the same way ID_0009 is populated: A1/B1 pull the
#    USUARIO/SENHA headers from ID_0001 via formula, C1/D1 are the category
#    headers, and row 2 is a sample data row.
$ws0010.Range("A1").Formula = "=ID_0001!A1"
$ws0010.Range("B1").Formula = "=ID_0001!B1"
$ws0010.Range("C1").Value = "CATEGORIA_1"
$ws0010.Range("D1").Value = "CATEGORIA_2"
$ws0010.Range("A1:D1").Font.Bold = $true

$ws0010.Range("A2").Value = "mmm"
$ws0010.Range("B2").Value = "mmm"
$ws0010.Range("C2").Value = "Monitors"
$ws0010.Range("D2").Value = "Laptops"

# 3. ID_0009's selection now spans both rows of its sample data
[void]$ws0009.Activate()
[void]$ws0009.Range("A1:D2").Select()

# 4. ID_0010 becomes the active / selected sheet and cell, finishing the
#    skeleton (this also clears tabSelected on ID_0008 automatically)
[void]$ws0010.Activate()
[void]$ws0010.Range("C2").Select()

Write-Output "ID_0010 scenario skeleton finished"
